$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44165
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("N2").Value = 20500
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20750
$ws.Range("Q2").Value = "$/caja 15 kilos"
$ws.Range("S2").Value = 1383
$ws.Range("T2").Value = 15
# Row 3
$ws.Range("D3").Value = 44165
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 17500
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 17750
$ws.Range("Q3").Value = "$/caja 15 kilos"
$ws.Range("S3").Value = 1183
$ws.Range("T3").Value = 15
# Row 4
$ws.Range("D4").Value = 44187
$ws.Range("K4").Value = "Dina"
$ws.Range("N4").Value = 22000
$ws.Range("O4").Value = 23000
$ws.Range("P4").Value = 22500
$ws.Range("Q4").Value = "$/caja 18 kilos"
$ws.Range("S4").Value = 1250
$ws.Range("T4").Value = 18
# Row 5
$ws.Range("D5").Value = 44186
$ws.Range("K5").Value = "Dina"
$ws.Range("L5").Value = "Especial"
$ws.Range("N5").Value = 22500
$ws.Range("O5").Value = 23000
$ws.Range("P5").Value = 22750
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("S5").Value = 1264
$ws.Range("T5").Value = 18
# Row 6
$ws.Range("D6").Value = 44168
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 23500
$ws.Range("O6").Value = 24000
$ws.Range("P6").Value = 23750
$ws.Range("Q6").Value = "$/caja 18 kilos"
$ws.Range("S6").Value = 1319
$ws.Range("T6").Value = 18
# Row 7
$ws.Range("D7").Value = 44162
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 20500
$ws.Range("O7").Value = 21000
$ws.Range("P7").Value = 20750
$ws.Range("S7").Value = 1383
# Row 8
$ws.Range("D8").Value = 44162
$ws.Range("K8").Value = "Castle Brite"
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 17500
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 17750
$ws.Range("S8").Value = 1183
# Row 9
$ws.Range("D9").Value = 44167
$ws.Range("K9").Value = "Castle Brite"
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 400
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 20500
$ws.Range("S9").Value = 1367
# Row 10
$ws.Range("D10").Value = 44167
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 360
$ws.Range("N10").Value = 17000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 17500
$ws.Range("S10").Value = 1167
# Row 11
$ws.Range("D11").Value = 44174
$ws.Range("M11").Value = 240
$ws.Range("N11").Value = 22500
$ws.Range("O11").Value = 23000
$ws.Range("P11").Value = 22750
$ws.Range("Q11").Value = "$/caja 18 kilos"
$ws.Range("S11").Value = 1264
$ws.Range("T11").Value = 18
# Row 14
$ws.Range("D14").Value = 44161
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 240
$ws.Range("N14").Value = 19500
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 19750
$ws.Range("Q14").Value = "$/caja 15 kilos"
$ws.Range("S14").Value = 1317
$ws.Range("T14").Value = 15
# Row 15
$ws.Range("D15").Value = 44161
$ws.Range("K15").Value = "Dina"
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 140
$ws.Range("N15").Value = 17500
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 17750
$ws.Range("Q15").Value = "$/caja 15 kilos"
$ws.Range("S15").Value = 1183
$ws.Range("T15").Value = 15
